# "MAR 12 EOD - After test w/ Richie"
#
# - add a new user: richard attfield, as user_01 (row 3 of user_list)
# - add a new response-tracking sheet "user_01" (a copy of the user_00
#   layout) with Richard's probe results
# - fix up a couple of values on user_00's "sect1 libB" section
# - leave the workbook with the user_list sheet active/selected

$wb = $excel.ActiveWorkbook

$wsList = $wb.Worksheets.Item("user_list")
$wsUser00 = $wb.Worksheets.Item("user_00")

# ---------------------------------------------------------------------
# 1. user_list: register the new user in row 3 (user_ID "01" is already
#    there; just the name was missing)
# ---------------------------------------------------------------------
$wsList.Range("B3").Value = "richard attfield (beta test)"

# ---------------------------------------------------------------------
# 2. user_00: a couple of probed values were corrected after the retest
# ---------------------------------------------------------------------
$wsUser00.Range("C6").Value = 1
$wsUser00.Range("D6").Value = 4
$wsUser00.Range("C8").Value = 1
$wsUser00.Range("D8").Value = 6

# ---------------------------------------------------------------------
# 3. Add the new "user_01" sheet (right after user_00) with Richard's
#    response data, laid out the same way as the other response sheets.
# ---------------------------------------------------------------------
$wsUser01 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wsUser00)
$wsUser01.Name = "user_01"

$wsUser01.Range("B1").Value = "Actual State"
$wsUser01.Range("C1").Value = "Probed State"
$wsUser01.Range("D1").Value = "Probed Confidence"

$wsUser01.Range("A2").Value = "sect1 libA"
$wsUser01.Range("B2").Value = 0
$wsUser01.Range("C2").Value = 0
$wsUser01.Range("D2").Value = 10
$wsUser01.Range("B3").Value = 1
$wsUser01.Range("C3").Value = 1
$wsUser01.Range("D3").Value = 10
$wsUser01.Range("B4").Value = 2
$wsUser01.Range("C4").Value = 2
$wsUser01.Range("D4").Value = 10

$wsUser01.Range("A6").Value = "sect1 libB"
$wsUser01.Range("B6").Value = 0
$wsUser01.Range("C6").Value = 0
$wsUser01.Range("D6").Value = 9
$wsUser01.Range("B7").Value = 1
$wsUser01.Range("C7").Value = 1
$wsUser01.Range("D7").Value = 8
$wsUser01.Range("B8").Value = 2
$wsUser01.Range("C8").Value = 2
$wsUser01.Range("D8").Value = 9

$wsUser01.Range("A10").Value = "sect3 libA"
$wsUser01.Range("B10").Value = 0
$wsUser01.Range("C10").Value = 0
$wsUser01.Range("D10").Value = 9
$wsUser01.Range("B11").Value = 1
$wsUser01.Range("C11").Value = 1
$wsUser01.Range("D11").Value = 9
$wsUser01.Range("B12").Value = 2
$wsUser01.Range("C12").Value = 2
$wsUser01.Range("D12").Value = 9

$wsUser01.Range("A14").Value = "sect3 libB"
$wsUser01.Range("B14").Value = 0
$wsUser01.Range("C14").Value = 0
$wsUser01.Range("D14").Value = 9
$wsUser01.Range("B15").Value = 1
$wsUser01.Range("C15").Value = 1
$wsUser01.Range("D15").Value = 9
$wsUser01.Range("B16").Value = 2
$wsUser01.Range("C16").Value = 2
$wsUser01.Range("D16").Value = 9

# leave user_01's own last selection on D8, like it was left in the
# original edit
$wsUser01.Activate()
$wsUser01.Range("D8").Select()

# ---------------------------------------------------------------------
# 4. Final view state: user_list becomes the active/selected tab again,
#    with B8 selected.
# ---------------------------------------------------------------------
$wsList.Activate()
$wsList.Range("B8").Select()
